# Scheduled runner update: refresh computed market-price columns (H..N)
# across the Behemoth_Profits leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2299.6875
$ws.Range("I33").Value = 2342.32
$ws.Range("J33").Value = 2147.4285
$ws.Range("K33").Value = 2342.32
$ws.Range("L33").Value = 2147.4285
$ws.Range("M33").Value = -2113.32
$ws.Range("N33").Value = -2605.4285

$ws.Range("H74").Value = 23529.4
$ws.Range("I74").Value = 28236.75
$ws.Range("K74").Value = 28236.75
$ws.Range("M74").Value = -27300.75

$ws.Range("H77").Value = 23529.4
$ws.Range("I77").Value = 28236.75
$ws.Range("K77").Value = 141183.75
$ws.Range("M77").Value = -136503.75

$ws.Range("H125").Value = 2059.3845
$ws.Range("I125").Value = 1720.2858
$ws.Range("J125").Value = 2455
$ws.Range("K125").Value = 15482.5722
$ws.Range("L125").Value = 22095
$ws.Range("M125").Value = -13022.5722
$ws.Range("N125").Value = -27015

$ws.Range("H132").Value = 2083.2258
$ws.Range("I132").Value = 1753.3077
$ws.Range("K132").Value = 5259.9231
$ws.Range("M132").Value = -2729.9231

$ws.Range("H135").Value = 19178.285
$ws.Range("I135").Value = 4616.6665
$ws.Range("K135").Value = 41549.9985
$ws.Range("M135").Value = -39014.9985

$ws.Range("H137").Value = 6132.25
$ws.Range("I137").Value = 2311.75
$ws.Range("J137").Value = 8042.5
$ws.Range("K137").Value = 6935.25
$ws.Range("L137").Value = 24127.5
$ws.Range("M137").Value = -4385.25
$ws.Range("N137").Value = -29227.5

$ws.Range("H138").Value = 1059645.5
$ws.Range("I138").Value = 676.5
$ws.Range("J138").Value = 1469569
$ws.Range("K138").Value = 2029.5
$ws.Range("L138").Value = 4408707
$ws.Range("M138").Value = 3110.5
$ws.Range("N138").Value = -4418987

$ws.Range("H141").Value = 2736.1304
$ws.Range("I141").Value = 2740.0476
$ws.Range("K141").Value = 8220.1428
$ws.Range("M141").Value = -3040.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 39923.832
$ws.Range("J31").Value = 110773
$ws.Range("L31").Value = 110773
$ws.Range("N31").Value = -111361

$ws.Range("H45").Value = 1760.7368
$ws.Range("I45").Value = 939.3333
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 939.3333
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -562.3333
$ws.Range("N45").Value = -3254

$ws.Range("H122").Value = 1600
$ws.Range("I122").Value = 1600
$ws.Range("K122").Value = 4800
$ws.Range("M122").Value = -2350

$ws.Range("H132").Value = 6593.84
$ws.Range("I132").Value = 3593.2354
$ws.Range("K132").Value = 10779.7062
$ws.Range("M132").Value = -8249.706200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2816.1667
$ws.Range("I86").Value = 2495.625
$ws.Range("J86").Value = 3457.25
$ws.Range("K86").Value = 2495.625
$ws.Range("L86").Value = 3457.25
$ws.Range("M86").Value = -1372.625
$ws.Range("N86").Value = -5703.25

$ws.Range("H89").Value = 2816.1667
$ws.Range("I89").Value = 2495.625
$ws.Range("J89").Value = 3457.25
$ws.Range("K89").Value = 12478.125
$ws.Range("L89").Value = 17286.25
$ws.Range("M89").Value = -6862.125
$ws.Range("N89").Value = -28518.25

$ws.Range("H94").Value = 1109.0526
$ws.Range("I94").Value = 1132.625
$ws.Range("J94").Value = 983.3333
$ws.Range("K94").Value = 1132.625
$ws.Range("L94").Value = 983.3333
$ws.Range("M94").Value = -681.625
$ws.Range("N94").Value = -1885.3333

$ws.Range("H96").Value = 43837.2
$ws.Range("J96").Value = 70874.60000000001
$ws.Range("L96").Value = 70874.60000000001
$ws.Range("N96").Value = -76366.60000000001

$ws.Range("H107").Value = 1872.6666
$ws.Range("I107").Value = 1646.7333
$ws.Range("K107").Value = 1646.7333
$ws.Range("M107").Value = 273.2666999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 870901.4
$ws.Range("I31").Value = 13990.833
$ws.Range("K31").Value = 13990.833
$ws.Range("M31").Value = -13695.833

$ws.Range("H34").Value = 870901.4
$ws.Range("I34").Value = 13990.833
$ws.Range("K34").Value = 13990.833
$ws.Range("M34").Value = -13788.833

$ws.Range("H58").Value = 2391.5386
$ws.Range("I58").Value = 2391.5386
$ws.Range("K58").Value = 2391.5386
$ws.Range("M58").Value = -2188.5386

$ws.Range("H105").Value = 1854.125
$ws.Range("I105").Value = 1923
$ws.Range("J105").Value = 1812.8
$ws.Range("K105").Value = 1923
$ws.Range("L105").Value = 1812.8
$ws.Range("M105").Value = -176
$ws.Range("N105").Value = -5306.8

$ws.Range("H134").Value = 479869.94
$ws.Range("I134").Value = 626579.25
$ws.Range("K134").Value = 1879737.75
$ws.Range("M134").Value = -1877202.75

$ws.Range("H136").Value = 2391.5386
$ws.Range("I136").Value = 2391.5386
$ws.Range("K136").Value = 7174.6158
$ws.Range("M136").Value = -4624.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 506578.5
$ws.Range("J39").Value = 749999.5
$ws.Range("L39").Value = 2249998.5
$ws.Range("N39").Value = -2250586.5

$ws.Range("H55").Value = 1000
$ws.Range("I55").Value = 1000
$ws.Range("K55").Value = 3000
$ws.Range("M55").Value = -2823

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6400.6
$ws.Range("I70").Value = 4002.6667
$ws.Range("J70").Value = 9997.5
$ws.Range("K70").Value = 4002.6667
$ws.Range("L70").Value = 9997.5
$ws.Range("M70").Value = -3732.6667
$ws.Range("N70").Value = -10537.5

$ws.Range("H73").Value = 6400.6
$ws.Range("I73").Value = 4002.6667
$ws.Range("J73").Value = 9997.5
$ws.Range("K73").Value = 4002.6667
$ws.Range("L73").Value = 9997.5
$ws.Range("M73").Value = -3066.6667
$ws.Range("N73").Value = -11869.5

$ws.Range("H102").Value = 2456.7917
$ws.Range("I102").Value = 2433.1738
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2433.1738
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -811.1738
$ws.Range("N102").Value = -6244

$ws.Range("H113").Value = 3917.7334
$ws.Range("I113").Value = 3755.077
$ws.Range("J113").Value = 4975
$ws.Range("K113").Value = 3755.077
$ws.Range("L113").Value = 4975
$ws.Range("M113").Value = -1585.077
$ws.Range("N113").Value = -9315

$ws.Range("H132").Value = 52633588
$ws.Range("I132").Value = 58825508
$ws.Range("K132").Value = 176476524
$ws.Range("M132").Value = -176473994

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 51546
$ws.Range("I7").Value = 3121
$ws.Range("K7").Value = 3121
$ws.Range("M7").Value = -3009

$ws.Range("H122").Value = 6844.4
$ws.Range("I122").Value = 5962.8887
$ws.Range("K122").Value = 17888.6661
$ws.Range("M122").Value = -15438.6661

$ws.Range("H126").Value = 51546
$ws.Range("I126").Value = 3121
$ws.Range("K126").Value = 9363
$ws.Range("M126").Value = -6893

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 22657.5
$ws.Range("J97").Value = 22657.5
$ws.Range("L97").Value = 22657.5
$ws.Range("N97").Value = -24639.5

$ws.Range("H119").Value = 65990
$ws.Range("J119").Value = 65990
$ws.Range("L119").Value = 65990
$ws.Range("N119").Value = -75666

$ws.Range("H122").Value = 5844.905
$ws.Range("J122").Value = 8443.25
$ws.Range("L122").Value = 25329.75
$ws.Range("N122").Value = -30229.75

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
